# Add an "Example" sheet to the NASA-TLX template workbook, update the
# directions text on the Template sheet, and rename the original sheet.

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet to "Template" -------------------------------
$template = $wb.Worksheets.Item(1)
$template.Name = "Template"

# --- Update the directions text to point at the new example sheet ---------
$template.Range("M3").Value = "Igonore all headings and align data with cell A1. Place weights directly next to scores. To record weights write the first word in all lower case letters in the correct cell. See example sheet."

# --- Add the new "Example" sheet right after "Template" --------------------
$example = $wb.Worksheets.Add($null, $template)
$example.Name = "Example"

# --- Fill in the example data ----------------------------------------------
$exampleData = @(
    @(45, 50, 55, 50, 5,  "effort",      "performance", "performance"),
    @(20, 15, 30, 55, 5,  "temporal",    "temporal",    "temporal"),
    @(10, 15, 25, 85, 5,  "temporal",    "temporal",    "effort"),
    @(5,  30, 65, 25, 10, "performance", "frustration", "frustration"),
    @(5,  65, 70, 45, 5,  "physical",    "temporal",    "physical"),
    @(35, 40, 35, 40, 15, "performance", "performance", "frustration")
)

for ($i = 0; $i -lt $exampleData.Length; $i++) {
    $row = $i + 1
    $values = $exampleData[$i]
    for ($col = 1; $col -le $values.Length; $col++) {
        $example.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

$labelRows = @(
    @("effort",      "effort",      "effort"),
    @("mental",      "mental",      "mental"),
    @("mental",      "effort",      "mental"),
    @("physical",    "physical",    "mental"),
    @("frustration", "frustration", "frustration"),
    @("mental",      "temporal",    "mental"),
    @("effort",      "temporal",    "temporal"),
    @("physical",    "physical",    "physical"),
    @("frustration", "mental",      "mental")
)

for ($i = 0; $i -lt $labelRows.Length; $i++) {
    $row = 7 + $i
    $values = $labelRows[$i]
    for ($col = 1; $col -le $values.Length; $col++) {
        $example.Cells.Item($row, 6 + $col - 1).Value = $values[$col - 1]
    }
}

# --- Put the focus back on the Template sheet (matches the saved state) ----
$template.Activate()
$template.Range("M4").Select() | Out-Null
